# Updated cryptos list - apply new prices/volumes and re-rank two coin pairs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# Row 2 - Bitcoin
Set-TextValue "D2" "65.694.56"
$ws.Range("E2").Value = "  +1.82%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.656.52"
$ws.Range("E3").Value = "  +0.87%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 - BNB
Set-TextValue "D5" "609.40"
$ws.Range("E5").Value = "  +2.48%  "

# Row 6 - Solana
Set-TextValue "D6" "156.80"
$ws.Range("E6").Value = "  +2.79%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.02%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.31%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "2.654.44"
$ws.Range("E9").Value = "  +0.88%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +8.75%  "

# Row 11 - now Toncoin (was Cardano)
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D11" "5.99"
$ws.Range("E11").Value = "  +3.00%  "

# Row 12 - now Cardano (was Toncoin)
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue "D12" "0.405"
$ws.Range("E12").Value = "  +2.02%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +1.53%  "

# Row 14 - Avalanche
Set-TextValue "D14" "30.08"
$ws.Range("E14").Value = "  +4.75%  "

# Row 15 - ShibaInu
Set-TextValue "D15" "0.0000203"
$ws.Range("E15").Value = "  +18.82%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue "D16" "3.137.12"
$ws.Range("E16").Value = "  +1.00%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "65.476.52"
$ws.Range("E17").Value = "  +1.66%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.649.62"
$ws.Range("E18").Value = "  -0.07%  "

# Row 19 - Chainlink
Set-TextValue "D19" "12.70"
$ws.Range("E19").Value = "  +3.51%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  +2.87%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "360.07"
$ws.Range("E21").Value = "  +2.48%  "

# Row 22 - Uniswap
Set-TextValue "D22" "7.45"
$ws.Range("E22").Value = "  +3.84%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.10%  "

# Row 24 - Litecoin
Set-TextValue "D24" "70.28"
$ws.Range("E24").Value = "  +4.10%  "

# Row 25 - SuiNetwork
Set-TextValue "D25" "1.70"
$ws.Range("E25").Value = "  -0.75%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextValue "D26" "9.60"
$ws.Range("E26").Value = "  +2.99%  "

# Row 27 - PEPE
Set-TextValue "D27" "0.0000107"
$ws.Range("E27").Value = "  +16.08%  "

# Row 28 - Fetch.AI
Set-TextValue "D28" "1.64"
$ws.Range("E28").Value = "  -0.95%  "

# Row 29 - Kaspa
$ws.Range("E29").Value = "  +2.53%  "

# Row 30 - Aptos
Set-TextValue "D30" "8.15"
$ws.Range("E30").Value = "  -1.83%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +6.44%  "

# Row 32 - Binance-PegBSC-USD
Set-TextValue "D32" "1.00"
$ws.Range("E32").Value = "  -2.37%  "

# Row 33 - Bittensor
Set-TextValue "D33" "534.34"
$ws.Range("E33").Value = "  -1.59%  "

# Row 34 - ImmutableX
Set-TextValue "D34" "1.80"
$ws.Range("E34").Value = "  -0.63%  "

# Row 35 - NEARProtocol
Set-TextValue "D35" "5.55"
$ws.Range("E35").Value = "  -1.31%  "

# Row 36 - RenderToken
$ws.Range("E36").Value = "  +3.27%  "

# Row 37 - PolygonEcosystemToken
Set-TextValue "D37" "0.433"
$ws.Range("E37").Value = "  +1.96%  "

# Row 38 - EthereumClassic
Set-TextValue "D38" "20.68"
$ws.Range("E38").Value = "  +2.53%  "

# Row 39 - Monero
Set-TextValue "D39" "163.45"
$ws.Range("E39").Value = "  -0.47%  "

# Row 40 - Stacks
Set-TextValue "D40" "2.00"
$ws.Range("E40").Value = "  -0.54%  "

# Row 41 - FirstDigitalUSD
Set-TextValue "D41" "1.00"
$ws.Range("E41").Value = "  -0.02%  "

# Row 42 - USDe
$ws.Range("E42").Value = "  +0.05%  "

# Row 43 - now OKB (was Aave)
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D43" "42.10"
$ws.Range("E43").Value = "  +1.19%  "

# Row 44 - now Aave (was OKB)
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D44" "165.62"
$ws.Range("E44").Value = "  -2.11%  "

# Row 45 - Filecoin
$ws.Range("E45").Value = "  +1.06%  "

# Row 46 - dogwifhat
Set-TextValue "D46" "2.34"
$ws.Range("E46").Value = "  +5.49%  "

# Row 47 - Hedera
Set-TextValue "D47" "0.0612"
$ws.Range("E47").Value = "  +2.96%  "

# Row 48 - InjectiveProtocol
Set-TextValue "D48" "23.11"
$ws.Range("E48").Value = "  -1.84%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  +4.95%  "

# Row 50 - Mantle
Set-TextValue "D50" "0.654"
$ws.Range("E50").Value = "  +1.59%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  +0.37%  "
